# Update latest output (run 165)
$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$schedule.Range("E2").Value = 331.5397897500001
$schedule.Range("F2").Value = 7.309078257275134
$schedule.Range("E3").Value = 426.124335
$schedule.Range("F3").Value = 28.18282638888889

# --- Detailed sheet updates ---
$detailed.Range("B17").Value = 0.51
$detailed.Range("B18").Value = -5.95032

$detailed.Range("B19").Value = 0.05995
$detailed.Range("C19").Value = "historical"

$detailed.Range("B20").Value = -5.2795
$detailed.Range("C20").Value = "historical"

$detailed.Range("B21").Value = -5.30722
$detailed.Range("B22").Value = -5.04585
$detailed.Range("B23").Value = -7.50157
$detailed.Range("B24").Value = -8.70809
$detailed.Range("B25").Value = -8
$detailed.Range("B26").Value = -12.01
$detailed.Range("B27").Value = -7.69975
$detailed.Range("B28").Value = -8.926069999999999
$detailed.Range("B29").Value = -7.20649
$detailed.Range("B30").Value = -5.50985
$detailed.Range("B31").Value = -5.71391
$detailed.Range("B32").Value = -12.01
$detailed.Range("B33").Value = -5.2795
$detailed.Range("B35").Value = -15.66234
$detailed.Range("B36").Value = -9.218389999999999
$detailed.Range("B37").Value = -7.97878
$detailed.Range("B38").Value = -0.46746
$detailed.Range("B39").Value = 7.26367
$detailed.Range("B40").Value = 36.25
$detailed.Range("B43").Value = 53.61259
$detailed.Range("B45").Value = 57.03541
$detailed.Range("B46").Value = 43.54764
$detailed.Range("B49").Value = 56.55141
